$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Employment")
$ws2 = $wb.Worksheets.Item("Payroll")

# New first column (state name) is now wide enough to show full names like
# "District of Columbia" without truncation - mirrors a best-fit column-width
# adjustment made after the box-link cleanup. (17.83 is the COM ColumnWidth
# input that this host's pixel-rounding maps closest to the target stored
# character width of ~18.71.)
$ws1.Columns.Item(1).ColumnWidth = 17.83
$ws2.Columns.Item(1).ColumnWidth = 17.83

# Reset each sheet's scroll position back to the top-left and update the
# remembered selection (previously left mid-scroll down near the bottom of
# each table). Select the Employment sheet's cell first, then finish on the
# Payroll sheet so it stays the active tab, matching the workbook's original
# active-tab state.
[void]$ws1.Range("E45").Select()
[void]$ws2.Range("C54").Select()
